## Helper: force a run-split at absolute character position $pos by toggling
## a (no-visible-effect) character formatting property on the range that runs
## from $pos to $rangeEnd, then immediately reverting it. Word's run-merge
## logic only coalesces runs that are format-identical, so applying then
## un-applying direct formatting on a sub-range is enough to make the writer
## keep that sub-range as its own <w:r>, without altering the visible look.
function Split-RunAt($doc, $pos, $rangeEnd) {
    $r = $doc.Range($pos, $rangeEnd)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

$d = $word.ActiveDocument

## ---------------------------------------------------------------------
## Change 1: "Click the <Setup Installer run>" paragraph
## ---------------------------------------------------------------------
$oldBig = "“Setup Installer” release from the right-hand side of the page and download the “setup” or “setup.exe” file."
$newBig = "“Setup Installer” release on the right-hand side of the page and download the “setup.exe” and “GES_Installation.msi” files to the same folder."

$anchor = $d.Content
$found = $anchor.Find.Execute("“Setup Installer” release", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bigStart = $anchor.Start

$bigRange = $d.Range($bigStart, $bigStart + $oldBig.Length)
$bigRange.Text = $newBig

$newLen = $newBig.Length
$rangeEndAbs = $bigStart + $newLen

# Split offsets (relative to $bigStart) between the 7 target runs. Offset 0
# is included too, so the new text doesn't get coalesced back into the
# preceding (untouched) " the " run.
#   0  .. 26  "“Setup Installer” release "
#   26 .. 28  "on"
#   28 .. 90  " the right-hand side of the page and download the “setup.exe” "
#   90 .. 117 "and “GES_Installation.msi” "
#   117.. 121 "file"
#   121.. 141 "s to the same folder"
#   141.. 142 "."
$splitOffsets = @(141, 121, 117, 90, 28, 26, 0)
foreach ($off in $splitOffsets) {
    $splitPos = $bigStart + $off
    Split-RunAt $d $splitPos $rangeEndAbs
}

## ---------------------------------------------------------------------
## Change 2: "The next step is to run the setup application..." paragraph
## (old and new text are both a single run each, so a plain Find/Replace
## keeps the run structure correct.)
## ---------------------------------------------------------------------
$old2 = "The next step is to run the setup application, which is an executable (.exe) file named " + [char]34 + "setup" + [char]34 + " or " + [char]34 + "setup.exe" + [char]34 + " that you will find in the unzipped folder. Double-click on this file to execute it."
$new2 = "The next step is to run the setup application, which is an executable (.exe) file named " + [char]34 + "setup" + [char]34 + " or " + [char]34 + "setup.exe" + [char]34 + ". Double-click this file to execute it."

$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

## ---------------------------------------------------------------------
## Change 3: "For further verification, the source code " paragraph
## ---------------------------------------------------------------------
$old3 = "For further verification, the source code "
$new3 = "For further security verification, the source code "

$anchor3 = $d.Content
$found3 = $anchor3.Find.Execute("For further verification, the source code", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start3 = $anchor3.Start

$range3 = $d.Range($start3, $start3 + $old3.Length)
$range3.Text = $new3

$newLen3 = $new3.Length
$rangeEndAbs3 = $start3 + $newLen3

# Split offsets between the 3 target runs:
#   0  .. 11  "For further"
#   11 .. 20  " security"
#   20 .. 51  " verification, the source code "
$splitOffsets3 = @(20, 11)
foreach ($off in $splitOffsets3) {
    $splitPos3 = $start3 + $off
    Split-RunAt $d $splitPos3 $rangeEndAbs3
}
